$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

# Update Runmode column (E) values from "Yes" to "no" for all rows except row 3
$ws.Range("E2").Value = "no"
$ws.Range("E4:E30").Value = "no"

# Update the view/selection so the active cell is E4 with the selection range E4:E30
$ws.Range("E4:E30").Select()
